$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Fix Implemented column (D2:D24) with detailed descriptions
$ws.Range("D2").Value = 'Added a session check to prevent logging in while already logged in'
$ws.Range("D3").Value = 'Enforced admin privilege checks before executing admin-specific transactions'
$ws.Range("D4").Value = 'Added a session check to prevent logout when no session is active'
$ws.Range("D5").Value = 'Implemented a check to ensure an account is associated with the session before withdrawing'
$ws.Range("D6").Value = 'Restricted withdrawal amount to a maximum of $500 using get_valid_number_input()'
$ws.Range("D7").Value = 'Added a loop that prevents withdrawals exceeding the available balance'
$ws.Range("D8").Value = 'Implemented checks for both sender and receiver accounts before proceeding'
$ws.Range("D9").Value = 'Ensured that the sender account matches the current session''s account'
$ws.Range("D10").Value = 'Restricted transfer amount to a maximum of $1000 using get_valid_number_input()'
$ws.Range("D11").Value = 'Added a loop that prevents transfers exceeding the sender''s balance'
$ws.Range("D12").Value = 'Implemented a check to ensure an account is associated with the session before paying bills'
$ws.Range("D13").Value = 'Restricted bill payment amount to a maximum of $2000 using get_valid_number_input()'
$ws.Range("D14").Value = 'Added a loop that prevents bill payments exceeding the account balance'
$ws.Range("D15").Value = 'Implemented a check to ensure an account is associated with the session before depositing'
$ws.Range("D16").Value = 'Enforced admin privilege checks before allowing account creation'
$ws.Range("D17").Value = 'Added a check to ensure account names do not exceed 20 characters'
$ws.Range("D18").Value = 'Restricted the initial balance input to a maximum of $99,999.99 using get_valid_number_input()'
$ws.Range("D19").Value = 'Enforced admin privilege checks before allowing account deletion'
$ws.Range("D20").Value = 'Implemented a check to ensure the account exists before deletion'
$ws.Range("D21").Value = 'Enforced admin privilege checks before allowing account disabling'
$ws.Range("D22").Value = 'Implemented a check to ensure the account exists before disabling'
$ws.Range("D23").Value = 'Enforced admin privilege checks before allowing plan changes'
$ws.Range("D24").Value = 'Implemented a check to ensure the account exists before changing the plan'

# Widen column D to fit the longer text
$ws.Columns.Item(4).ColumnWidth = 78.3

# Update the active selection
$ws.Range("D25").Select() | Out-Null

